$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tyrese Haliburton"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Indiana Pacers"

$ws.Range("A5").Value = "Darius Garland"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Cleveland Cavaliers"

$ws.Range("A6").Value = "Keegan Murray"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Sacramento Kings"

$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"

$ws.Range("A14").Value = "OG Anunoby"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "New York Knicks"

$ws.Range("A16").Value = "Stephen Curry"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Golden State Warriors"

$ws.Range("A19").Value = "Daniel Gafford"
$ws.Range("B19").Value = "PF,C"
$ws.Range("C19").Value = "Dallas Mavericks"

